$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column F (shifts old F,G,H -> G,H,I)
$ws.Columns("F").Insert()

# Header for new column
$ws.Range("F8").Value = "URL Change and Final changes"

# Populate new F column with "Done" for rows that had an F-value originally (now in G)
$ws.Range("F9").Value = "Done"
$ws.Range("F10").Value = "Done"
$ws.Range("F11").Value = "Done"
$ws.Range("F13").Value = "Done"
$ws.Range("F14").Value = "Done"
$ws.Range("F15").Value = "Done"
$ws.Range("F16").Value = "Done"
$ws.Range("F17").Value = "Done"
$ws.Range("F18").Value = "Done"
$ws.Range("F19").Value = "Done"
$ws.Range("F20").Value = "Done"
$ws.Range("F21").Value = "Done"

# Update the view
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F25").Select()
